# Generate Report for Handoff
#
# A new handoff Xliff generation event occurred for the
# "4bbc1c0c-778b-4977-8f48-9a20cca500de.md" file (row 7 on every sheet),
# so the recorded "Latest Handoff" timestamps for that row need to be
# updated on the Overview sheet as well as on each locale sheet.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-28 10:50:31"

# zh-cn sheet: column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-28 10:50:26"

# de-de sheet: column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-28 10:50:31"
